# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header (A1): refresh time 19:52 -> 20:22
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 20:22"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 603488
$ws.Range("C4").Value = 16547
$ws.Range("D4").Value = 38131
$ws.Range("E4").Value = 540163
$ws.Range("F4").Value = 12828
$ws.Range("G4").Value = 1554
$ws.Range("H4").Value = 25194

# Brasil (row 17)
$ws.Range("B17").Value = 24232
$ws.Range("C17").Value = 802
$ws.Range("E17").Value = 19808

# Austria (row 20)
$ws.Range("B20").Value = 14224
$ws.Range("C20").Value = 183
$ws.Range("E20").Value = 6207

# Israel (row 21)
$ws.Range("B21").Value = 12046
$ws.Range("C21").Value = 460
$ws.Range("D21").Value = 2195
$ws.Range("E21").Value = 9728
$ws.Range("F21").Value = 175
$ws.Range("G21").Value = 7
$ws.Range("H21").Value = 123

# Ecuador (row 29)
$ws.Range("D29").Value = 696
$ws.Range("E29").Value = 6552

# Row 41: country label swaps from "Indonesia" to "Emiratos Arabes Unidos"
$ws.Range("A41").Value = "Emiratos Arabes Unidos"
$ws.Range("B41").Value = 4933
$ws.Range("C41").Value = 412
$ws.Range("D41").Value = 933
$ws.Range("E41").Value = 3972
$ws.Range("F41").Value = 1
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 28

# Row 42: country label swaps from "Emiratos Arabes Unidos" to "Indonesia"
$ws.Range("A42").Value = "Indonesia"
$ws.Range("B42").Value = 4839
$ws.Range("C42").Value = 282
$ws.Range("D42").Value = 426
$ws.Range("E42").Value = 3954
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 459
$ws.Range("H42").Value = 25

# Luxemburgo (row 47)
$ws.Range("B47").Value = 3307
$ws.Range("C47").Value = 15
$ws.Range("E47").Value = 2738

# Georgia (row 112)
$ws.Range("D112").Value = 69
$ws.Range("E112").Value = 224

# Isla de Man (row 115)
$ws.Range("B115").Value = 254
$ws.Range("C115").Value = 12
$ws.Range("E115").Value = 111
